$wb = $excel.ActiveWorkbook

# Sheet1: update A1 and A2 values
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("A1").Value = "29038,98"
$ws1.Range("A2").Value = "82,52"

# data sheet: change header labels and clear stale cells
$ws2 = $wb.Worksheets.Item("data")
$ws2.Range("A1").Value = "-"
$ws2.Range("B1").Value = "-"
$ws2.Range("A2").Value = ""
$ws2.Range("B2").Value = ""
$ws2.Range("A3").Value = ""
$ws2.Range("B3").Value = ""
